$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 36: module code changes from MEEN30020 to MEEN30140 (new row added for BIOMED path)
$ws.Range("A36").Value = "MEEN30140"

# Row 37: module code changes from MEEN20020 to MEEN30030,
# Level (column D) changes from 2 to 3, Trimester (column H) changes from Aut to Spr
$ws.Range("A37").Value = "MEEN30030"
$ws.Range("D37").Value = 3
$ws.Range("H37").Value = "Spr"

# Scroll the view down so row 16 is the top visible row
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# Select B37:H37 (active cell B37), matching the final on-screen selection
$ws.Range("B37:H37").Select()
